$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ligand/receptor/edge TPM-derived statistics (columns G-J, M-P, Q-T)
# for rows 2-10 with newly recomputed TPM values.

$ws.Range("G2").Value = 72.92148999999999
$ws.Range("H2").Value = 218.76447
$ws.Range("I2").Value = 0.2015977907456805
$ws.Range("J2").Value = 0.2015977907456805
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 147.8202705085133
$ws.Range("R2").Value = 1330.38243457662
$ws.Range("S2").Value = 0.001329796395518224
$ws.Range("T2").Value = 0.001329796395518224
$ws.Range("G3").Value = 72.92148999999999
$ws.Range("H3").Value = 218.76447
$ws.Range("I3").Value = 0.2015977907456805
$ws.Range("J3").Value = 0.2015977907456805
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 18700.30084692161
$ws.Range("R3").Value = 168302.7076222945
$ws.Range("S3").Value = 0.1682285695716581
$ws.Range("T3").Value = 0.168228569571658
$ws.Range("G4").Value = 72.92148999999999
$ws.Range("H4").Value = 218.76447
$ws.Range("I4").Value = 0.2015977907456805
$ws.Range("J4").Value = 0.2015977907456805
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 3561.504944409177
$ws.Range("R4").Value = 32053.54449968259
$ws.Range("S4").Value = 0.03203942477850421
$ws.Range("T4").Value = 0.0320394247785042
$ws.Range("I5").Value = 0.591090693015494
$ws.Range("J5").Value = 0.591090693015494
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 433.4134109973477
$ws.Range("R5").Value = 3900.72069897613
$ws.Range("S5").Value = 0.003899002415100696
$ws.Range("T5").Value = 0.003899002415100696
$ws.Range("I6").Value = 0.591090693015494
$ws.Range("J6").Value = 0.591090693015494
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.4932511482656078
$ws.Range("T6").Value = 0.4932511482656077
$ws.Range("I7").Value = 0.591090693015494
$ws.Range("J7").Value = 0.591090693015494
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 10442.43797505026
$ws.Range("R7").Value = 93981.94177545229
$ws.Range("S7").Value = 0.09394054233478559
$ws.Range("T7").Value = 0.09394054233478555
$ws.Range("G8").Value = 74.98824566666667
$ws.Range("H8").Value = 224.964737
$ws.Range("I8").Value = 0.2073115162388255
$ws.Range("J8").Value = 0.2073115162388255
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 152.0098226106669
$ws.Range("R8").Value = 1368.088403496002
$ws.Range("S8").Value = 0.001367485754799696
$ws.Range("T8").Value = 0.001367485754799695
$ws.Range("G9").Value = 74.98824566666667
$ws.Range("H9").Value = 224.964737
$ws.Range("I9").Value = 0.2073115162388255
$ws.Range("J9").Value = 0.2073115162388255
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 19230.30856815368
$ws.Range("R9").Value = 173072.7771133831
$ws.Range("S9").Value = 0.1729965378270716
$ws.Range("T9").Value = 0.1729965378270716
$ws.Range("G10").Value = 74.98824566666667
$ws.Range("H10").Value = 224.964737
$ws.Range("I10").Value = 0.2073115162388255
$ws.Range("J10").Value = 0.2073115162388255
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 3662.445840237266
$ws.Range("R10").Value = 32962.01256213539
$ws.Range("S10").Value = 0.03294749265695424
$ws.Range("T10").Value = 0.03294749265695422
